# Instruction List.xlsx - LWL, LWR, mfthilo instr tests
#
# This script applies the changes observed between the "before" and "after"
# snapshots of "Instruction List.xlsx":
#   * Several "Tested"/"Implemented" (columns N / L) flags flip from
#     "No" -> "Yes" for BNE, MFHI, MTHI, MFLO, MTLO, LWL and LWR rows.
#   * Row 28 (BEQ)'s "Tested" flag also flips No -> Yes.
#   * The stray review comments in column M (rows 27 & 28) are removed -
#     row 27's M cell becomes fully empty, row 28's M cell becomes an
#     empty-but-styled (underlined font) cell.
#   * The two now-unused shared strings ("Check for all branches negative
#     condition" and "signextend(imm<<2) is bad") fall out of the shared
#     string table as a natural consequence of removing their only
#     references.
#   * The view scrolls down to row 37 and the selection moves to N47,
#     with the zoom normalized to 100%.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column M cleanup (rows 27 & 28) ------------------------------------

# Row 27: comment "Check for all branches negative condition" is deleted
# outright - the cell goes away entirely.
$ws.Range("M27").ClearContents()

# Row 28: comment "signextend(imm<<2) is bad" is deleted too, but the cell
# keeps an (empty) underlined style that was left behind.
$ws.Range("M28").ClearContents()
$ws.Range("M28").Font.Underline = $true

# --- Tested / Implemented flags: No -> Yes ------------------------------

$ws.Range("N28").Value = "Yes"   # BEQ     - Tested
$ws.Range("N31").Value = "Yes"   # BNE     - Tested
$ws.Range("N43").Value = "Yes"   # MFHI    - Tested
$ws.Range("N44").Value = "Yes"   # MTHI    - Tested
$ws.Range("N45").Value = "Yes"   # MFLO    - Tested
$ws.Range("N46").Value = "Yes"   # MTLO    - Tested
$ws.Range("L54").Value = "Yes"   # LWL     - Implemented
$ws.Range("L55").Value = "Yes"   # LWR     - Implemented

# --- View state: scroll position, selection, zoom -----------------------

$ws.Range("N47").Select()

try { $excel.ActiveWindow.ScrollRow = 37 } catch {}
try { $excel.ActiveWindow.ScrollColumn = 1 } catch {}
try { $excel.ActiveWindow.Zoom = 100 } catch {}
